$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (column D) and volume/1h (column E) values
# scraped on 2023-01-21, preserving original text formatting (inline string cells).

$c = $ws.Range("D2")
$c.Value = "'304.49"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'2.41%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'35.67"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "'12.64%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.094"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'2.20%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.07803"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'1.29%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'2.266"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'1.30%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'8.113"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'2.75%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'4.036"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'6.58%"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'0.28%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.09577"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'-3.40%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.1827"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'4.48%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.08559"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'2.07%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.03427"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'5.76%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.09949"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'1.25%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.001481"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'0.86%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'0.005736"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'-0.01%"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'-1.21%"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'2.179"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'-0.79%"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'2.92%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.1322"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'0.49%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'4.555"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'12.05%"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'-1.33%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.04681"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'3.70%"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'2.69%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.004542"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'4.08%"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'1.04%"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'-19.55%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.01770"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'4.10%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.04716"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'1.82%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.007948"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'5.63%"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'2.08%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.008015"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'-17.59%"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'11.33%"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'-6.26%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00006200"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'2.75%"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.00000000751"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'1.05%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'4.053"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'45.07%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'0.002694"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'36.25%"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'1.05%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.0002003"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'1.05%"
$c.Style = "Normal"
